# language.xlsx edit script
# - shorten all cycle duration to 50 sec., more tweaks to levels, adjusted gardener AI, some language tweak
#
# Concretely (as captured by the canonical OOXML diff):
#   1. Bump several "VoiceDuration" values in column C (rows 104-114) on both
#      the "en" and "es" sheets.
#   2. Remove the (now-unused) "intro_climate_polar_01" / "Brrrtz! ..." pair
#      of rows (row 201 on both sheets) which collapses the shared-string
#      table and reflows every row below it.
#   3. Switch the active sheet/tab to "es" and update the remembered
#      selection on each sheet.

$wb = $excel.ActiveWorkbook
$wsEn = $wb.Worksheets.Item("en")
$wsEs = $wb.Worksheets.Item("es")

# --- 1. VoiceDuration tweaks (column C), identical on both language sheets ---
foreach ($ws in @($wsEn, $wsEs)) {
    $ws.Range("C104").Value = 1.5
    $ws.Range("C106").Value = 1.5
    $ws.Range("C108").Value = 1.5
    $ws.Range("C110").Value = 2
    $ws.Range("C112").Value = 1.5
    $ws.Range("C114").Value = 1.5
}

# --- 2. Drop the obsolete "intro_climate_polar_01" / "Brrrtz!..." row ---
# Row 201 holds Key="intro_climate_polar_01" / Value="Brrrtz! This climate is
# as cold as it can get!" on "en", and the matching (untranslated) Key on
# "es". Deleting the whole row on both sheets removes both now-orphaned
# shared strings and shifts every following row up by one.
$wsEn.Rows.Item(201).EntireRow.Delete()
$wsEs.Rows.Item(201).EntireRow.Delete()

# --- 3. View state: select "es" as the active tab, refresh remembered selections ---
$wsEn.Activate()
$excel.ActiveWindow.ScrollRow = 164
$excel.ActiveWindow.ScrollColumn = 1
$wsEn.Range("B199").Select()

$wsEs.Activate()
$excel.ActiveWindow.ScrollRow = 133
$excel.ActiveWindow.ScrollColumn = 1
$wsEs.Range("B181").Select()
